$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '53.132.16'
$ws.Range("E2").Value = '  -12.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.372.50'
$ws.Range("E3").Value = '  -17.70%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '455.29'
$ws.Range("E5").Value = '  -12.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.70'
$ws.Range("E6").Value = '  -8.77%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.35%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.479'
$ws.Range("E8").Value = '  -10.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.381.02'
$ws.Range("E9").Value = '  -17.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0928'
$ws.Range("E10").Value = '  -12.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.31'
$ws.Range("E11").Value = '  -13.18%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.314'
$ws.Range("E12").Value = '  -11.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.121'
$ws.Range("E13").Value = '  -5.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.774.93'
$ws.Range("E14").Value = '  -18.21%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '53.105.66'
$ws.Range("E15").Value = '  -12.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.35'
$ws.Range("E16").Value = '  -12.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000124'
$ws.Range("E17").Value = '  -11.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.362.12'
$ws.Range("E18").Value = '  -18.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.14'
$ws.Range("E19").Value = '  -14.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '304.82'
$ws.Range("E20").Value = '  -12.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.15'
$ws.Range("E21").Value = '  -19.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.02'
$ws.Range("E22").Value = '  +2.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.33'
$ws.Range("E24").Value = '  -17.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '55.01'
$ws.Range("E25").Value = '  -14.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.996'
$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.379'
$ws.Range("E27").Value = '  -14.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.419.82'
$ws.Range("E28").Value = '  -19.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.147'
$ws.Range("E29").Value = '  -16.51%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.08'
$ws.Range("E30").Value = '  -7.77%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.995'
$ws.Range("E31").Value = '  -0.47%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0₃0690'
$ws.Range("E32").Value = '  -18.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '146.16'
$ws.Range("E33").Value = '  -4.77%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.40'
$ws.Range("E34").Value = '  -9.96%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.38'
$ws.Range("E35").Value = '  -16.82%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.93'
$ws.Range("E36").Value = '  -10.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.829'
$ws.Range("E37").Value = '  -14.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.45'
$ws.Range("E38").Value = '  -20.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.03'
$ws.Range("E39").Value = '  -12.56%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.994'
$ws.Range("E40").Value = '  -0.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '32.86'
$ws.Range("E41").Value = '  -11.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.590'
$ws.Range("E42").Value = '  -8.70%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.25'
$ws.Range("E43").Value = '  -10.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0519'
$ws.Range("E44").Value = '  -9.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.13'
$ws.Range("E45").Value = '  -2.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.951.08'
$ws.Range("E46").Value = '  -13.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.20'
$ws.Range("E47").Value = '  -16.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0213'
$ws.Range("E48").Value = '  -8.90%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0852'
$ws.Range("E49").Value = '  -5.73%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.18'
$ws.Range("E50").Value = '  -13.07%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.40'
$ws.Range("E51").Value = '  -18.20%  '

